# Update ids for Hypertension workflow.
#
# The first worksheet ("Initial Diagnosis (ID)") is renamed to
# "Hypertension Initial Dx (H)" and every "ID-*" endpoint id living in
# column C of that sheet is renamed to the matching "H-*" id. The
# previously selected cell on that sheet (C9) is also reset to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab.
$ws.Name = "Hypertension Initial Dx (H)"

# Rename each ID-* endpoint identifier to its H-* counterpart. The order
# below matches the order the ids were first introduced so the shared
# string table ends up appended in the same sequence.
$ws.Range("C3").Value2 = "H-ExcludedUnder18"
$ws.Range("C4").Value2 = "H-ExcludedPregnant"
$ws.Range("C5").Value2 = "H-ExcludedEndStageRenalDisease"
$ws.Range("C6").Value2 = "H-ExcludedNormalBP"
$ws.Range("C7").Value2 = "H-HypertensiveEmergencySBP"
$ws.Range("C8").Value2 = "H-HypertensiveEmergencyDBP"
$ws.Range("C9").Value2 = "H-MonitoringPreexistingHTN"
$ws.Range("C10").Value2 = "H-RecommendMoreBPs"
$ws.Range("C11").Value2 = "H-HTNStage2LastBP"
$ws.Range("C12").Value2 = "H-HTNStage2AverageBP"
$ws.Range("C13").Value2 = "H-ConsiderHTNStage2"
$ws.Range("C14").Value2 = "H-ConsiderHTNStage1"
$ws.Range("C15").Value2 = "H-PrescribeHBPABPMonitoring"
$ws.Range("C16").Value2 = "H-PrescribeAmbulatoryBPMonitoring"
$ws.Range("C17").Value2 = "H-NoFurtherAction"
$ws.Range("C2").Value2 = "H-ExcludedOver80"

# Reset the saved selection on the sheet from C9 to C2.
$ws.Activate()
$ws.Range("C2").Select()
